$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set B6 formula to produce a #NUM! error (overflow from repeated multiplication)
$ws.Range("B6").Formula = "=1E+99*1E+99*1E+99*1E+99"

# Update the active selection to E14
$ws.Range("E14").Select()
